$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates: 종가(D), RSI(E), 5일수익률(F)
$ws.Range("D2").Value = 91589.32000000001
$ws.Range("E2").Value = 47.5
$ws.Range("F2").Value = 0.33

# MACRO_SCORE (column N) updated for all data rows 2-6
$ws.Range("N2").Value = 66.09241856096124
$ws.Range("N3").Value = 66.09241856096124
$ws.Range("N4").Value = 66.09241856096124
$ws.Range("N5").Value = 66.09241856096124
$ws.Range("N6").Value = 66.09241856096124
